# 남근우 일정표.xlsx — weekly log update
# - Add a "5th week" status entry to Sheet2 (new rows of text)
# - Make Sheet2 the active/selected sheet (was Sheet1)

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = "5주차 한 일"
$ws2.Range("A2").Value = "DB 수정: 맵 사이즈의 최대 저장 크기 증가"

$ws2.Activate()
$ws2.Range("A3").Select() | Out-Null
